$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(13).Insert()

$ws.Cells.Item(13, 1).Value = 5
$ws.Cells.Item(13, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(13, 3).Value = "Maule"
$ws.Cells.Item(13, 4).Value = 44462
$ws.Cells.Item(13, 5).Value = 7
$ws.Cells.Item(13, 6).Value = 100112031
$ws.Cells.Item(13, 7).Value = "Poroto verde"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 150
$ws.Cells.Item(13, 11).Value = 30000
$ws.Cells.Item(13, 12).Value = 30000
$ws.Cells.Item(13, 13).Value = 30000
$ws.Cells.Item(13, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(13, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(13, 16).Value = 1200
$ws.Cells.Item(13, 17).Value = 25
$ws.Cells.Item(13, 18).Value = "Hortaliza"
